$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = '20.838.52'
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").Value = '1.524.96'
$ws.Range("E3").Value = '  +4.55%  '
Set-TextValue $ws.Range("D4") '1.006'
$ws.Range("E4").Value = '  -0.21%  '
Set-TextValue $ws.Range("D5") '0.9583'
$ws.Range("E5").Value = '  +1.54%  '
Set-TextValue $ws.Range("D6") '280.08'
$ws.Range("E6").Value = '  +2.08%  '
Set-TextValue $ws.Range("D7") '0.3582'
$ws.Range("E7").Value = '  -1.72%  '
Set-TextValue $ws.Range("D8") '0.3135'
$ws.Range("E8").Value = '  +2.28%  '
Set-TextValue $ws.Range("D9") '1.100'
$ws.Range("E9").Value = '  +6.71%  '
Set-TextValue $ws.Range("D10") '40.10'
$ws.Range("E10").Value = '  +0.87%  '
Set-TextValue $ws.Range("D11") '0.06729'
$ws.Range("E11").Value = '  +2.76%  '
Set-TextValue $ws.Range("D12") '1.001'
$ws.Range("E12").Value = '  +0.23%  '
Set-TextValue $ws.Range("D13") '18.57'
$ws.Range("E13").Value = '  +4.33%  '
Set-TextValue $ws.Range("D14") '5.594'
$ws.Range("E14").Value = '  +3.93%  '
Set-TextValue $ws.Range("D15") '6.275'
$ws.Range("E15").Value = '  +2.91%  '
Set-TextValue $ws.Range("D16") '0.9589'
$ws.Range("E16").Value = '  -0.10%  '
Set-TextValue $ws.Range("D17") '0.00001032'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").Value = '1.516.95'
$ws.Range("E18").Value = '  +4.32%  '
Set-TextValue $ws.Range("D19") '0.06028'
$ws.Range("E19").Value = '  +4.90%  '
Set-TextValue $ws.Range("D20") '70.73'
$ws.Range("E20").Value = '  +1.90%  '
Set-TextValue $ws.Range("D21") '5.612'
$ws.Range("E21").Value = '  +3.63%  '
Set-TextValue $ws.Range("D22") '14.84'
$ws.Range("E22").Value = '  +3.09%  '
Set-TextValue $ws.Range("D23") '11.32'
$ws.Range("E23").Value = '  +4.45%  '
Set-TextValue $ws.Range("D24") '2.297'
$ws.Range("E24").Value = '  +2.81%  '
$ws.Range("D25").Value = '20.881.79'
$ws.Range("E25").Value = '  +2.40%  '
Set-TextValue $ws.Range("D26") '147.03'
$ws.Range("E26").Value = '  +4.75%  '
Set-TextValue $ws.Range("D27") '2.163'
$ws.Range("E27").Value = '  +4.21%  '
Set-TextValue $ws.Range("D28") '17.43'
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("D29").Value = '1.683.44'
$ws.Range("E29").Value = '  +4.50%  '
Set-TextValue $ws.Range("D30") '116.64'
$ws.Range("E30").Value = '  +4.58%  '
Set-TextValue $ws.Range("D31") '4.052'
$ws.Range("E31").Value = '  +7.83%  '
Set-TextValue $ws.Range("D32") '5.111'
$ws.Range("E32").Value = '  +5.41%  '
Set-TextValue $ws.Range("D33") '0.8349'
$ws.Range("E33").Value = '  +6.41%  '
Set-TextValue $ws.Range("D34") '0.07987'
$ws.Range("E34").Value = '  +2.58%  '
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D35") '1.477'
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D36") '1.207'
$ws.Range("E36").Value = '  +7.11%  '
Set-TextValue $ws.Range("D37") '4.854'
$ws.Range("E37").Value = '  +4.48%  '
Set-TextValue $ws.Range("D38") '0.05798'
$ws.Range("E38").Value = '  +1.62%  '
Set-TextValue $ws.Range("D39") '0.02067'
$ws.Range("E39").Value = '  +2.44%  '
Set-TextValue $ws.Range("D40") '10.53'
$ws.Range("E40").Value = '  +2.06%  '
Set-TextValue $ws.Range("D41") '0.9592'
$ws.Range("E41").Value = '  +0.86%  '
Set-TextValue $ws.Range("D42") '0.1884'
$ws.Range("E42").Value = '  +1.38%  '
Set-TextValue $ws.Range("D43") '7.528'
$ws.Range("E43").Value = '  +1.87%  '
Set-TextValue $ws.Range("D44") '0.5343'
$ws.Range("E44").Value = '  +1.85%  '
Set-TextValue $ws.Range("D45") '3.548'
$ws.Range("E45").Value = '  +1.94%  '
Set-TextValue $ws.Range("D46") '12.29'
$ws.Range("E46").Value = '  +3.59%  '
Set-TextValue $ws.Range("D47") '120.67'
$ws.Range("E47").Value = '  +3.48%  '
Set-TextValue $ws.Range("D48") '0.5341'
$ws.Range("E48").Value = '  +4.09%  '
Set-TextValue $ws.Range("D49") '1.845'
$ws.Range("E49").Value = '  +5.92%  '
Set-TextValue $ws.Range("D50") '0.06520'
$ws.Range("E50").Value = '  +1.75%  '
Set-TextValue $ws.Range("D51") '0.9881'
$ws.Range("E51").Value = '  +0.04%  '
